# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# per the commit diff (crypto price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.601.15"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.906.56"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.340"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "2.180.00"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.12%  "
$ws.Range("D14").Value = "1.949.51"
$ws.Range("E14").Value = "  +5.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.696"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "35.564.36"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "72.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +30.55%  "
$ws.Range("E27").Value = "  +8.32%  "
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.976"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +31.22%  "
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0566"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("E35").Value = "  +8.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.27%  "
$ws.Range("E39").Value = "  +5.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "1.360.71"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0599"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.16%  "
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +29.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +40.17%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("D50").Value = "2.090.43"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("E51").Value = "  +3.60%  "
